# Update the "Working with ROS" slide: turn the trailing "provided tutorial."
# words of the *Note paragraph into a hyperlink that points at the tutorial
# file, matching the target OOXML (split run + a:hlinkClick + endParaRPr).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shape = $s.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

# The note is the third paragraph of the content placeholder.
$paraCount = $textRange.Paragraphs().Count
$notePara = $textRange.Paragraphs($paraCount)

$linkText = "provided tutorial."
$fullText = $notePara.Text
$startIndex = $fullText.IndexOf($linkText)

# PowerPoint character ranges are 1-based.
$startPos = $startIndex + 1
$length = $linkText.Length

$linkRange = $notePara.Characters($startPos, $length)
$hyperlink = $linkRange.ActionSettings.Item(1).Hyperlink
$hyperlink.Address = "ROS Tutorial.docx"
